$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1092.2354
$ws.Range("I98").Value = 726.5
$ws.Range("J98").Value = 1970
$ws.Range("K98").Value = 726.5
$ws.Range("L98").Value = 1970
$ws.Range("M98").Value = 771.5
$ws.Range("N98").Value = -4966
$ws.Range("H107").Value = 13889729
$ws.Range("I107").Value = 14706477
$ws.Range("K107").Value = 14706477
$ws.Range("M107").Value = -14704557
$ws.Range("H116").Value = 6586.6924
$ws.Range("I116").Value = 7132.7827
$ws.Range("J116").Value = 2400
$ws.Range("K116").Value = 7132.7827
$ws.Range("L116").Value = 2400
$ws.Range("M116").Value = -3690.7827
$ws.Range("N116").Value = -9284
$ws.Range("H121").Value = 1202.9333
$ws.Range("J121").Value = 1241
$ws.Range("L121").Value = 3723
$ws.Range("N121").Value = -7217
$ws.Range("H122").Value = 1092.2354
$ws.Range("I122").Value = 726.5
$ws.Range("J122").Value = 1970
$ws.Range("K122").Value = 2179.5
$ws.Range("L122").Value = 5910
$ws.Range("M122").Value = 270.5
$ws.Range("N122").Value = -10810
$ws.Range("H131").Value = 3926.3333
$ws.Range("I131").Value = 1674.375
$ws.Range("J131").Value = 6500
$ws.Range("K131").Value = 5023.125
$ws.Range("L131").Value = 19500
$ws.Range("M131").Value = 16.875
$ws.Range("N131").Value = -29580
$ws.Range("H133").Value = 50000
$ws.Range("J133").Value = 50000
$ws.Range("L133").Value = 50000
$ws.Range("N133").Value = -60120
$ws.Range("H137").Value = 1745
$ws.Range("I137").Value = 1421.8387
$ws.Range("K137").Value = 4265.5161
$ws.Range("M137").Value = -1715.5161

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 162379.9
$ws.Range("I61").Value = 3902.7908
$ws.Range("J61").Value = 503105.7
$ws.Range("K61").Value = 3902.7908
$ws.Range("L61").Value = 503105.7
$ws.Range("M61").Value = -3690.7908
$ws.Range("N61").Value = -503529.7
$ws.Range("H123").Value = 38929
$ws.Range("J123").Value = 38929
$ws.Range("L123").Value = 38929
$ws.Range("N123").Value = -48729
$ws.Range("H136").Value = 162379.9
$ws.Range("I136").Value = 3902.7908
$ws.Range("J136").Value = 503105.7
$ws.Range("K136").Value = 11708.3724
$ws.Range("L136").Value = 1509317.1
$ws.Range("M136").Value = -9158.3724
$ws.Range("N136").Value = -1514417.1

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H60").Value = 0
$ws.Range("J60").Value = 0
$ws.Range("L60").ClearContents()
$ws.Range("N60").Value = 0
$ws.Range("H64").Value = 259.58334
$ws.Range("I64").Value = 100
$ws.Range("J64").Value = 291.5
$ws.Range("K64").Value = 100
$ws.Range("L64").Value = 291.5
$ws.Range("M64").Value = 125
$ws.Range("N64").Value = -741.5
$ws.Range("H67").Value = 259.58334
$ws.Range("I67").Value = 100
$ws.Range("J67").Value = 291.5
$ws.Range("K67").Value = 100
$ws.Range("L67").Value = 291.5
$ws.Range("M67").Value = 680
$ws.Range("N67").Value = -1851.5
$ws.Range("H134").Value = 23937.844
$ws.Range("I134").Value = 4895.6
$ws.Range("J134").Value = 93182.37
$ws.Range("K134").Value = 14686.8
$ws.Range("L134").Value = 279547.11
$ws.Range("M134").Value = -12151.8
$ws.Range("N134").Value = -284617.11

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1467.5834
$ws.Range("I16").Value = 1283
$ws.Range("J16").Value = 2113.625
$ws.Range("K16").Value = 1283
$ws.Range("L16").Value = 2113.625
$ws.Range("M16").Value = -996
$ws.Range("N16").Value = -2687.625
$ws.Range("H31").Value = 10110356
$ws.Range("I31").Value = 1668.5358
$ws.Range("J31").Value = 22416584
$ws.Range("K31").Value = 1668.5358
$ws.Range("L31").Value = 22416584
$ws.Range("M31").Value = -1373.5358
$ws.Range("N31").Value = -22417174
$ws.Range("H34").Value = 10110356
$ws.Range("I34").Value = 1668.5358
$ws.Range("J34").Value = 22416584
$ws.Range("K34").Value = 1668.5358
$ws.Range("L34").Value = 22416584
$ws.Range("M34").Value = -1466.5358
$ws.Range("N34").Value = -22416988
$ws.Range("H113").Value = 1467.5834
$ws.Range("I113").Value = 1283
$ws.Range("J113").Value = 2113.625
$ws.Range("K113").Value = 1283
$ws.Range("L113").Value = 2113.625
$ws.Range("M113").Value = 887
$ws.Range("N113").Value = -6453.625
$ws.Range("H132").Value = 5407778
$ws.Range("I132").Value = 8334850.5
$ws.Range("J132").Value = 3952
$ws.Range("K132").Value = 25004551.5
$ws.Range("L132").Value = 11856
$ws.Range("M132").Value = -25002021.5
$ws.Range("N132").Value = -16916

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 2602.7778
$ws.Range("J34").Value = 2883.3333
$ws.Range("L34").Value = 8649.999899999999
$ws.Range("N34").Value = -8817.999899999999
$ws.Range("H39").Value = 2400
$ws.Range("J39").Value = 2400
$ws.Range("L39").Value = 7200
$ws.Range("N39").Value = -7788
$ws.Range("H55").Value = 7071.4287
$ws.Range("J55").Value = 7071.4287
$ws.Range("L55").Value = 21214.2861
$ws.Range("N55").Value = -21568.2861
$ws.Range("H113").Value = 1053106.6
$ws.Range("I113").Value = 1316253.6
$ws.Range("J113").Value = 526812.8
$ws.Range("K113").Value = 3948760.8
$ws.Range("L113").Value = 1580438.4
$ws.Range("M113").Value = -3946590.8
$ws.Range("N113").Value = -1584778.4
$ws.Range("H131").Value = 3573056.5
$ws.Range("I131").Value = 100000000
$ws.Range("J131").Value = 1688.1852
$ws.Range("K131").Value = 300000000
$ws.Range("L131").Value = 5064.5556
$ws.Range("M131").Value = -299994960
$ws.Range("N131").Value = -15144.5556

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 104
$ws.Range("I3").Value = 104
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 104
$ws.Range("L3").Value = 0
$ws.Range("M3").ClearContents()
$ws.Range("N3").Value = 12
$ws.Range("H10").Value = 3983333.2
$ws.Range("I10").Value = 5500000
$ws.Range("J10").Value = 950000
$ws.Range("K10").Value = 5500000
$ws.Range("L10").Value = 950000
$ws.Range("M10").Value = -5499831
$ws.Range("N10").Value = -950338
$ws.Range("H11").Value = 21571428
$ws.Range("I11").Value = 21200000
$ws.Range("K11").Value = 21200000
$ws.Range("M11").Value = -21199861
$ws.Range("H12").Value = 5147500
$ws.Range("I12").Value = 5147500
$ws.Range("K12").Value = 5147500
$ws.Range("M12").Value = -5147360
$ws.Range("H102").Value = 2721.7407
$ws.Range("I102").Value = 2590.318
$ws.Range("K102").Value = 2590.318
$ws.Range("M102").Value = -968.3180000000002
$ws.Range("H132").Value = 2418052.8
$ws.Range("I132").Value = 3548341
$ws.Range("J132").Value = 3346.5454
$ws.Range("K132").Value = 10645023
$ws.Range("L132").Value = 10039.6362
$ws.Range("M132").Value = -10642493
$ws.Range("N132").Value = -15099.6362
$ws.Range("H140").Value = 26267.4
$ws.Range("J140").Value = 26267.4
$ws.Range("L140").Value = 26267.4
$ws.Range("N140").Value = -36627.4
$ws.Range("H141").Value = 65398.8
$ws.Range("J141").Value = 65398.8
$ws.Range("L141").Value = 65398.8
$ws.Range("N141").Value = -75758.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").ClearContents()
$ws.Range("I40").Value = 111113750
$ws.Range("J40").Value = 3102.5
$ws.Range("K40").Value = 111113750
$ws.Range("L40").Value = 3102.5
$ws.Range("M40").Value = -111113614
$ws.Range("N40").Value = -3374.5
$ws.Range("H61").Value = 1993.8125
$ws.Range("I61").Value = 1564.3572
$ws.Range("J61").Value = 5000
$ws.Range("K61").Value = 1564.3572
$ws.Range("L61").Value = 5000
$ws.Range("M61").Value = -1362.3572
$ws.Range("N61").Value = -5404
$ws.Range("H113").Value = 1993.8125
$ws.Range("I113").Value = 1564.3572
$ws.Range("J113").Value = 5000
$ws.Range("K113").Value = 1564.3572
$ws.Range("L113").Value = 5000
$ws.Range("M113").Value = 605.6428000000001
$ws.Range("N113").Value = -9340

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 69502
$ws.Range("J2").Value = 69502
$ws.Range("L2").Value = 69502
$ws.Range("N2").Value = -69726
$ws.Range("H122").Value = 618.94116
$ws.Range("I122").Value = 549
$ws.Range("J122").Value = 846.25
$ws.Range("K122").Value = 1647
$ws.Range("L122").Value = 2538.75
$ws.Range("M122").Value = 803
$ws.Range("N122").Value = -7438.75
$ws.Range("H136").Value = 2648198.2
$ws.Range("I136").Value = 2562.7576
$ws.Range("J136").Value = 5558397
$ws.Range("K136").Value = 7688.2728
$ws.Range("L136").Value = 16675191
$ws.Range("M136").Value = -5138.2728
$ws.Range("N136").Value = -16680291
$ws.Range("H141").Value = 62266.273
$ws.Range("J141").Value = 62266.273
$ws.Range("L141").Value = 62266.273
$ws.Range("N141").Value = -72626.273
